$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily COVID-19 data rows for Slovenia (2020-05-25 .. 2020-05-28),
# appended below the existing table data (which ended at row 75 / 2020-05-24).
$newRows = @(
    @(43976, 75770, 754, 1469, 0, 9, 2, 6, 108, 1),
    @(43977, 76579, 809, 1471, 2, 8, 2, 2, 108, 0),
    @(43978, 77210, 631, 1473, 2, 7, 2, 1, 108, 0),
    @(43979, 77916, 706, 1473, 0, 7, 2, 0, 108, 0)
)

$firstNewRow = 76
$lastNewRow = 79

# Match formatting used by the preceding data rows before filling in values,
# so the appended cells keep the table's date / number look instead of
# inheriting the plain column defaults.
$rDate = $ws.Range("A${firstNewRow}:A${lastNewRow}")
$rDate.NumberFormat = "d/\ m/\ yyyy;@"
$rDate.Font.Name = "Calibri Light"
$rDate.Font.Size = 10
$rDate.HorizontalAlignment = -4152
$rDate.VerticalAlignment = -4160

$rTestedAll = $ws.Range("B${firstNewRow}:B${lastNewRow}")
$rTestedAll.NumberFormat = "#,##0"
$rTestedAll.Font.Name = "Calibri Light"
$rTestedAll.Font.Size = 10
$rTestedAll.HorizontalAlignment = -4152

$rRest = $ws.Range("C${firstNewRow}:J${lastNewRow}")
$rRest.Font.Name = "Calibri Light"
$rRest.Font.Size = 10
$rRest.HorizontalAlignment = -4152

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $firstNewRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}

# Grow the table (ListObject) so it covers the newly added rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J$lastNewRow"))

# Match the saved selection state from the edit (last appended row selected).
$ws.Range("A$lastNewRow`:J$lastNewRow").Select()
